$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.75   # Current Capital
$summary.Range("B4").Value = 0.74      # Total P&L $
$summary.Range("B5").Value = 0.51      # Total P&L %
$summary.Range("B6").Value = 29        # Total Trades
$summary.Range("B7").Value = 12        # Winning Trades
$summary.Range("B9").Value = 41.38     # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.75     # Capital
$status.Range("D4").Value = 29         # Trades
$status.Range("E4").Value = 0.74       # P&L $
$status.Range("F4").Value = 0.75       # P&L %
$status.Range("G4").Value = 41.38      # Win Rate %

# --- All Trades sheet (Trade #29 is in row 30) ---
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G30").Value = 0.75          # Exit Price
$allTrades.Range("H30").Value = "CLOSED"      # Status
$allTrades.Range("I30").Value = 13.6364       # P&L %
$allTrades.Range("J30").Value = 0.09          # P&L $
$allTrades.Range("K30").Value = 100.75        # Capital After
$allTrades.Range("P30").Value = "early_exit"  # Exit Reason
$allTrades.Range("Q30").Value = 0.13          # Duration (min)

# --- MarketMaking sheet (same trade mirrored, row 30) ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G30").Value = 0.75          # Exit Price
$marketMaking.Range("H30").Value = "CLOSED"      # Status
$marketMaking.Range("I30").Value = 13.6364       # P&L %
$marketMaking.Range("J30").Value = 0.09          # P&L $
$marketMaking.Range("K30").Value = 100.75        # Capital After
$marketMaking.Range("P30").Value = "early_exit"  # Exit Reason
$marketMaking.Range("Q30").Value = 0.13          # Duration (min)
